# Attendance Tracker — update roster data, add a hyperlinked-email column,
# append a new attendee row, and touch up sheet selections / column widths
# to match the author's final save state.

$wb = $excel.ActiveWorkbook

# Helper: write a value to a cell as literal TEXT (never let Excel's
# auto-detection reinterpret a date-looking or number-looking string as a
# real date/number). We flip the cell to Text format, assign the value,
# then paste-special just the *formats* from a pristine, never-touched
# cell on the same sheet so the cell doesn't end up carrying a stray
# explicit style index in the saved file.
function Set-TextValue {
    param(
        $WorkSheet,
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $cell = $WorkSheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $blank = $WorkSheet.Cells.Item(500, 200)
    $blank.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------
# Sheet: Person_Master
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Person_Master")

# Jiraiya -> Jiraiya Ogata
$ws1.Range("B7").Value = "Jiraiya Ogata"

# Row 13: Hiei -> Orochimaru (id 13 -> 14, status change, follow-up date,
# now baptized, new mobile number)
$ws1.Range("A13").Value = 14
$ws1.Range("B13").Value = "Orochimaru"
Set-TextValue $ws1 13 3 "2026-02-07"
$ws1.Range("D13").Value = "For follow-up"
$ws1.Range("G13").Value = "Yes"
Set-TextValue $ws1 13 14 "1234567800"

# ---------------------------------------------------------------------
# Sheet: Attendance_Table
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Attendance_Table")

# Jiraiya -> Jiraiya Ogata (two historical attendance rows)
$ws2.Range("B7").Value = "Jiraiya Ogata"
$ws2.Range("B8").Value = "Jiraiya Ogata"

# Row 16 gets reassigned to Hinata Hyuuga's Feb-07 check-in...
$ws2.Range("A16").Value = 10
$ws2.Range("B16").Value = "Hinata Hyuuga"
Set-TextValue $ws2 16 4 "2026-02-07"

# ...and Orochimaru's new attendance becomes row 17
$ws2.Range("A17").Value = 14
$ws2.Range("B17").Value = "Orochimaru"
$ws2.Range("C17").Value = 2026
Set-TextValue $ws2 17 4 "2026-02-07"

# ---------------------------------------------------------------------
# Sheet: Cell_Group_Master
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Cell_Group_Master")

$ws3.Range("B2").Value = "Madara Uchiha"
$ws3.Range("C2").Value = "madara@uchiha.com"
$ws3.Hyperlinks.Add($ws3.Range("C2"), "mailto:madara@uchiha.com") | Out-Null

$ws3.Range("B3").Value = "Ippo Makunouchi"
$ws3.Range("C3").Value = "dempseyroll@hni.com"
$ws3.Hyperlinks.Add($ws3.Range("C3"), "mailto:dempseyroll@hni.com") | Out-Null

$ws3.Range("B4").Value = "Steve Armstrong"
$ws3.Range("C4").Value = "steve.armstrong@voltesv.com"
$ws3.Hyperlinks.Add($ws3.Range("C4"), "mailto:steve.armstrong@voltesv.com") | Out-Null

$ws3.Range("B5").Value = "Kakashi Hatake"
$ws3.Range("C5").Value = "thesixth@hokage.com"
$ws3.Hyperlinks.Add($ws3.Range("C5"), "mailto:thesixth@hokage.com") | Out-Null

$ws3.Columns.Item(3).ColumnWidth = 24.25

$ws3.Activate()
$ws3.Range("A7").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet: Ministry_Master
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Ministry_Master")

$ws4.Range("B2").Value = "Robert Akizuki"
$ws4.Range("B3").Value = "Jamie Robinson"
$ws4.Range("B4").Value = "Orihime Inoue"
$ws4.Range("B5").Value = "Frieren Tribbiani"
$ws4.Range("B6").Value = "Nico Robyn"
$ws4.Range("B7").Value = "Dai Sawamura"
$ws4.Range("B8").Value = "Yoruichi Shihouin"

$ws4.Columns.Item(1).ColumnWidth = 20.1
$ws4.Columns.Item(2).ColumnWidth = 15.1

$ws4.Activate()
$ws4.Range("B9").Select() | Out-Null
